$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the values that were in row 3 and row 4 ---
$ws.Range("A3").Value = 33
$ws.Range("B3").Value = "T931100609002"
$ws.Range("A4").Value = 24
$ws.Range("B4").Value = "F931100609041"

# --- Step 2: extend formatting (border / alignment style used by column A)
#     down through the newly added rows 5-31, matching row 3/4's look. ---
$ws.Range("A3:B3").Copy()
$ws.Range("A5:B31").PasteSpecial(-4122)   # xlPasteFormats

# --- Step 3: fill in the newly appended rows' data (rows 5 through 31) ---
$newData = @(
    @(30, "V931325309014"),
    @(36, "W931101108060"),
    @(3,  "Q931325208064"),
    @(4,  "H931325209012"),
    @(25, "M931252916068"),
    @(29, "Q931321008053"),
    @(10, "B931412016036"),
    @(11, "G931101109060"),
    @(26, "Y931321110015"),
    @(0,  "Y888201710013"),
    @(32, "K931101109004"),
    @(35, "D931100609028"),
    @(31, "X886463320016"),
    @(23, "Q931100609020"),
    @(15, "M931252710007"),
    @(19, "C931321610014"),
    @(16, "P931383310002"),
    @(13, "U931412020025"),
    @(12, "X931325210023"),
    @(9,  "W931101109061"),
    @(8,  "N931100609007"),
    @(7,  "Q879418719002"),
    @(6,  "D931100608056"),
    @(5,  "A931383810034"),
    @(2,  "R928218115049"),
    @(20, "U931101109019"),
    @(37, "T931100609029")
)

$startRow = 5
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
}
